$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the two trailing "section header + data" rows down by one row to make
# room for the new "amr_surveillance_annexC" requirement row at row 14.
# Using Range.Copy(Destination) (rather than Rows.Insert) avoids introducing
# spurious new style entries in the workbook.

# old row 15 (data_indicators_report) -> row 16
$ws.Range("A15:D15").Copy($ws.Range("A16:D16"))
# old row 14 (Generating Supplementary report header) -> row 15
$ws.Range("A14:D14").Copy($ws.Range("A15:D15"))

# Build the new row 14 by cloning the formatting/pattern of row 13
# (amr_surveillance_annexB), which also gives D14 the correct explanatory
# text (shared with D6:D13) without needing to re-read it.
$ws.Range("A13:D13").Copy($ws.Range("A14:D14"))

$ws.Range("A14").Value = "amr_surveillance_annexC"
$ws.Range("B14").Value = "yes"
$ws.Range("C14").Value = "Required"
